$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three oldest years (2007-2009) which live in rows 2-4; this
# shifts 2010..2020 (previously rows 5-15) up to rows 2-12.
$ws.Range("A2:J4").Delete(-4162) | Out-Null  # xlShiftUp

# Copy the year-label cell formatting (bold/centered/bordered style) from
# the row above down onto the new row so the new year label matches the
# look of every other row in column A.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Append 2021年, whose detailed breakdown (columns B:I) is not yet
# published -- only the grand total in column J is known, matching how
# 2019年/2020年 are already represented (blank measure cells, populated
# total).
$ws.Cells.Item(13, 1).Value = "2021年"
for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item(13, $c).Value = " "
}
$ws.Cells.Item(13, 10).Value = 8850
